# Apply the weekly cryptos-list refresh (prices / 1h % changes / two row swaps)
# scraped by the GitHub Actions job into Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "76.578.60"

# Row 3
$ws.Range("D3").Value = "3.032.80"
$ws.Range("E3").Value = "  +3.87%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D5").Value = "201.38"
$ws.Range("E5").Value = "  +0.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D6").Value = "634.76"
$ws.Range("E6").Value = "  +5.88%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D9").Value = "0.202"
$ws.Range("E9").Value = "  +2.44%  "

# Row 10
$ws.Range("D10").Value = "3.026.05"
$ws.Range("E10").Value = "  +3.72%  "

# Row 11
$ws.Range("E11").Value = "  +1.00%  "

# Row 12
$ws.Range("E12").Value = "  -0.10%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D13").Value = "5.04"
$ws.Range("E13").Value = "  +2.97%  "

# Row 14
$ws.Range("D14").Value = "3.569.04"
$ws.Range("E14").Value = "  +3.25%  "

# Row 15
$ws.Range("E15").Value = "  +6.49%  "

# Row 16
$ws.Range("D16").Value = "76.395.52"
$ws.Range("E16").Value = "  +0.67%  "

# Row 17
$ws.Range("E17").Value = "  -0.98%  "

# Row 18
$ws.Range("D18").Value = "3.007.04"
$ws.Range("E18").Value = "  +3.09%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D19").Value = "13.62"
$ws.Range("E19").Value = "  +6.46%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D20").Value = "9.01"
$ws.Range("E20").Value = "  +1.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D21").Value = "374.89"
$ws.Range("E21").Value = "  -0.90%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D22").Value = "2.28"
$ws.Range("E22").Value = "  -2.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D23").Value = "4.31"
$ws.Range("E23").Value = "  +2.71%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D24").Value = "72.91"
$ws.Range("E24").Value = "  +2.05%  "

# Row 25
$ws.Range("E25").Value = "  +2.15%  "

# Row 26
$ws.Range("E26").Value = "  +0.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D27").Value = "4.38"
$ws.Range("E27").Value = "  +3.31%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D28").Value = "9.83"
$ws.Range("E28").Value = "  +0.87%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D29").Value = "0.0000108"
$ws.Range("E29").Value = "  -1.04%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.51%  "

# Row 31
$ws.Range("E31").Value = "  +7.39%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  -0.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D33").Value = "513.67"
$ws.Range("E33").Value = "  +1.46%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D34").Value = "1.96"
$ws.Range("E34").Value = "  +8.22%  "

# Row 35
$ws.Range("E35").Value = "  -0.05%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D36").Value = "20.60"
$ws.Range("E36").Value = "  +2.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D37").Value = "163.19"
$ws.Range("E37").Value = "  -1.05%  "

# Row 38
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D38").Value = "0.385"
$ws.Range("E38").Value = "  +12.14%  "

# Row 39
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D39").Value = "20.01"
$ws.Range("E39").Value = "  +1.51%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D40").Value = "188.17"
$ws.Range("E40").Value = "  +4.44%  "

# Row 41
$ws.Range("E41").Value = "  +9.36%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D42").Value = "0.112"
$ws.Range("E42").Value = "  -1.35%  "

# Row 43
$ws.Range("E43").Value = "  +0.52%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D44").Value = "5.00"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D45").Value = "42.64"
$ws.Range("E45").Value = "  +6.08%  "

# Row 46
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D46").Value = "1.24"
$ws.Range("E46").Value = "  +3.02%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D47").Value = "1.65"
$ws.Range("E47").Value = "  -1.15%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D48").Value = "0.606"
$ws.Range("E48").Value = "  +4.95%  "

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D49").Value = "0.712"
$ws.Range("E49").Value = "  +7.48%  "

# Row 50
$ws.Range("E50").Value = "  -1.13%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"   # keep numeric-looking price as text
$ws.Range("D51").Value = "3.86"
$ws.Range("E51").Value = "  +3.55%  "
